$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-09-10 Tuesday"; new = "2024-09-11 Wednesday"},
    @{old = "430×6=2580"; new = "648×5=3240"},
    @{old = "891×9=8019"; new = "780×2=1560"},
    @{old = "813×3=2439"; new = "279×3=837"},
    @{old = "961×7=6727"; new = "661×5=3305"},
    @{old = "878×8=7024"; new = "815×6=4890"},
    @{old = "302×2=604"; new = "382×6=2292"},
    @{old = "109×9=981"; new = "142×8=1136"},
    @{old = "783×3=2349"; new = "390×5=1950"},
    @{old = "689×9=6201"; new = "326×3=978"},
    @{old = "582×4=2328"; new = "545×9=4905"},
    @{old = "905×8=7240"; new = "735×3=2205"},
    @{old = "174×2=348"; new = "816×6=4896"},
    @{old = "293×7=2051"; new = "277×8=2216"},
    @{old = "970×8=7760"; new = "972×4=3888"},
    @{old = "408×3=1224"; new = "296×9=2664"},
    @{old = "504×9=4536"; new = "354×2=708"},
    @{old = "657×4=2628"; new = "764×2=1528"},
    @{old = "992×9=8928"; new = "120×6=720"},
    @{old = "349×3=1047"; new = "171×9=1539"},
    @{old = "798×7=5586"; new = "607×3=1821"},
    @{old = "430×3=1290"; new = "962×8=7696"},
    @{old = "742×7=5194"; new = "842×3=2526"},
    @{old = "422×5=2110"; new = "736×9=6624"},
    @{old = "153×3=459"; new = "772×5=3860"},
    @{old = "157×3=471"; new = "624×9=5616"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
